$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

$ws.Range("AW2").Value = 156.272523
$ws.Range("AK3").Value = 47.110833
$ws.Range("AQ4").Value = 19.061875
$ws.Range("AQ5").Value = 0.230058
$ws.Range("AW6").Value = 111.218183
$ws.Range("AW7").Value = 114.004734
$ws.Range("AK8").Value = 14.10985
$ws.Range("AW9").Value = 149.042234
$ws.Range("AW10").Value = 98.08287
$ws.Range("AW11").Value = 82.166157
$ws.Range("AW12").Value = 83.06853
$ws.Range("AN13").Value = 25.651933
$ws.Range("AW14").Value = 161.126227
$ws.Range("AW15").Value = 146.941377
$ws.Range("AW16").Value = 116.069907
$ws.Range("AW17").Value = 118.367211
$ws.Range("AN18").Value = 21.240243
$ws.Range("AQ19").Value = 8.803680999999999
$ws.Range("AW20").Value = 85.95403899999999
$ws.Range("AW21").Value = 78.042801
$ws.Range("AK22").Value = 14.109861
$ws.Range("AK23").Value = 14.109907
$ws.Range("AQ24").Value = 0.229931
$ws.Range("AW25").Value = 111.21765
$ws.Range("AQ26").Value = 0.061759
$ws.Range("AK27").Value = 14.110012
$ws.Range("AW28").Value = 146.946887
$ws.Range("AW29").Value = 68.03037
$ws.Range("AW30").Value = 77.12385399999999
$ws.Range("AW31").Value = 7.2639
$ws.Range("AW32").Value = 149.037141
$ws.Range("AW33").Value = 85.265856
$ws.Range("AW34").Value = 6.045301
$ws.Range("AQ35").Value = 0.0625
$ws.Range("AW36").Value = 113.311655
$ws.Range("AW37").Value = 121.097396
$ws.Range("AW38").Value = 8.329294000000001
$ws.Range("AK39").Value = 14.109884
$ws.Range("AK40").Value = 14.10978
$ws.Range("AK41").Value = 14.110093
$ws.Range("AK42").Value = 14.109803
$ws.Range("AQ43").Value = 19.061748
$ws.Range("AW44").Value = 70.348657
$ws.Range("AW45").Value = 85.265243
$ws.Range("AW46").Value = 36.382014
$ws.Range("AW47").Value = 36.380799
$ws.Range("AK48").Value = 36.159398
$ws.Range("AW49").Value = 113.312639
$ws.Range("AW50").Value = 114.006146
$ws.Range("AW51").Value = 21.148762
$ws.Range("AW52").Value = 127.221632
$ws.Range("AW53").Value = 104.979167
$ws.Range("AW54").Value = 152.083819
$ws.Range("AW55").Value = 142.132639
$ws.Range("AK56").Value = 36.959641
$ws.Range("AW57").Value = 114.005313
$ws.Range("AW58").Value = 8.333403000000001
$ws.Range("AQ59").Value = 19.061262
$ws.Range("AW60").Value = 50.006887
$ws.Range("AW61").Value = 117.260995
$ws.Range("AW62").Value = 78.062141
$ws.Range("AK63").Value = 9.135612999999999
$ws.Range("AQ64").Value = 1.328403
$ws.Range("AW65").Value = 153.376956
$ws.Range("AW66").Value = 156.26985
$ws.Range("AW67").Value = 114.021088
$ws.Range("AW68").Value = 71.12261599999999
$ws.Range("AW69").Value = 14.155174
$ws.Range("AK70").Value = 14.110069
$ws.Range("AK71").Value = 14.109826
$ws.Range("AW72").Value = 7.159213
$ws.Range("AK73").Value = 36.159479
$ws.Range("AW74").Value = 128.092523
$ws.Range("AW75").Value = 65.11873799999999
$ws.Range("AW76").Value = 0.045729
$ws.Range("AW77").Value = 69.327384
$ws.Range("AN78").Value = 7.275382
$ws.Range("AW79").Value = 142.13603
$ws.Range("AW80").Value = 36.381794
$ws.Range("AW81").Value = 21.297153
$ws.Range("AW82").Value = 15.99191
$ws.Range("AW83").Value = 6.043958
$ws.Range("AW84").Value = 121.161354
$ws.Range("AW85").Value = 152.0361
$ws.Range("AW86").Value = 100.005764
$ws.Range("AW87").Value = 69.325671
$ws.Range("AW88").Value = 49.999734
$ws.Range("AQ89").Value = 34.091539
$ws.Range("AW90").Value = 106.008275
$ws.Range("AW91").Value = 146.940995
$ws.Range("AW92").Value = 56.381424
$ws.Range("AK93").Value = 9.133171000000001
$ws.Range("AW94").Value = 114.004537
$ws.Range("AK95").Value = 14.109734
$ws.Range("AW96").Value = 153.37735
$ws.Range("AW97").Value = 64.992014
$ws.Range("AQ98").Value = 19.061979
$ws.Range("AW99").Value = 83.06119200000001
$ws.Range("AW100").Value = 75.137049
$ws.Range("AW101").Value = 149.126956
$ws.Range("AW102").Value = 105.123947
$ws.Range("AW103").Value = 44.382697
$ws.Range("AK104").Value = 47.110799
$ws.Range("AK105").Value = 9.135636999999999
$ws.Range("AQ106").Value = 2.287778
$ws.Range("AW107").Value = 103.560822
$ws.Range("AW108").Value = 112.20434
$ws.Range("AW109").Value = 65.047963
$ws.Range("AK110").Value = 14.109954
$ws.Range("AN111").Value = 7.275336
$ws.Range("AW112").Value = 14.136655
$ws.Range("AK113").Value = 14.109769
$ws.Range("AW114").Value = 156.272813
$ws.Range("AW115").Value = 142.083692
$ws.Range("AW116").Value = 152.034236
$ws.Range("AW117").Value = 44.3811
$ws.Range("AW118").Value = 15.982546
$ws.Range("AQ119").Value = 0.06433999999999999
$ws.Range("AW120").Value = 131.002905
$ws.Range("AW121").Value = 105.143958
$ws.Range("AW122").Value = 83.014398
$ws.Range("AW123").Value = 69.323519
$ws.Range("AK124").Value = 14.109896
$ws.Range("AK125").Value = 36.959745
$ws.Range("AW126").Value = 146.941609
$ws.Range("AK127").Value = 14.109803
$ws.Range("AN128").Value = 7.275405
$ws.Range("AW129").Value = 142.036204
$ws.Range("AK130").Value = 36.159514
$ws.Range("AK131").Value = 14.10978
$ws.Range("AW132").Value = 56.243773
$ws.Range("AW133").Value = 50.395185
$ws.Range("AW134").Value = 160.086227
$ws.Range("AW135").Value = 85.926586
$ws.Range("AW136").Value = 44.380556
$ws.Range("AK137").Value = 14.109873
$ws.Range("AN138").Value = 21.165116
$ws.Range("AW139").Value = 142.086898
$ws.Range("AW140").Value = 146.943264
$ws.Range("AW141").Value = 83.12388900000001
$ws.Range("AK142").Value = 14.109815
$ws.Range("AW143").Value = 117.260567
$ws.Range("AW144").Value = 114.005035
$ws.Range("AW145").Value = 77.12488399999999
$ws.Range("AK146").Value = 14.110058
$ws.Range("AW147").Value = 49.25706
$ws.Range("AW148").Value = 49.057118
$ws.Range("AQ149").Value = 34.091424
$ws.Range("AW150").Value = 70.95525499999999
$ws.Range("AW151").Value = 61.020486
$ws.Range("AW152").Value = 91.09066
$ws.Range("AW153").Value = 44.379861
$ws.Range("AK154").Value = 14.109919
$ws.Range("AQ155").Value = 26.071539
$ws.Range("AW156").Value = 113.312442
$ws.Range("AW157").Value = 69.329132
$ws.Range("AK158").Value = 9.133171000000001
$ws.Range("AN159").Value = 2.142639
$ws.Range("AW160").Value = 103.335255
$ws.Range("AW161").Value = 19.081701
$ws.Range("AK162").Value = 14.110058
$ws.Range("AW163").Value = 152.037037
$ws.Range("AW164").Value = 128.387338
$ws.Range("AW165").Value = 68.030162
$ws.Range("AK166").Value = 14.109745
$ws.Range("AW167").Value = 103.340197
$ws.Range("AW168").Value = 156.273206
$ws.Range("AW169").Value = 78.041273
$ws.Range("AW170").Value = 149.036586
$ws.Range("AW171").Value = 140.143426
$ws.Range("AW172").Value = 78.06202500000001
$ws.Range("AW173").Value = 7.159051
$ws.Range("AK174").Value = 36.159387
$ws.Range("AW175").Value = 149.036227
$ws.Range("AW176").Value = 105.139769
$ws.Range("AW177").Value = 50.394468
$ws.Range("AW178").Value = 82.954178
$ws.Range("AN179").Value = 6.107118
$ws.Range("AW180").Value = 156.280775
$ws.Range("AW181").Value = 91.102847
$ws.Range("AW182").Value = 69.326644
$ws.Range("AW183").Value = 50.398009
$ws.Range("AW184").Value = 105.143588
$ws.Range("AK185").Value = 14.109896
$ws.Range("AK186").Value = 36.159514
$ws.Range("AW187").Value = 149.127234
$ws.Range("AW188").Value = 163.348067
$ws.Range("AW189").Value = 69.332488
$ws.Range("AW190").Value = 43.119884
$ws.Range("AW191").Value = 136.363773
$ws.Range("AW192").Value = 146.942176
$ws.Range("AW193").Value = 75.04572899999999
$ws.Range("AW194").Value = 130.999745
$ws.Range("AW195").Value = 78.969167
$ws.Range("AK196").Value = 9.135624999999999
$ws.Range("AW197").Value = 152.01037
$ws.Range("AW198").Value = 149.12669
$ws.Range("AW199").Value = 148.048912
$ws.Range("AW200").Value = 160.031493
$ws.Range("AW201").Value = 162.171887
$ws.Range("AW202").Value = 152.085069
$ws.Range("AW203").Value = 114.006759
$ws.Range("AW204").Value = 19.349769
$ws.Range("AK205").Value = 36.159398
$ws.Range("AW206").Value = 114.024734
$ws.Range("AW207").Value = 161.127211
$ws.Range("AW208").Value = 152.026134
$ws.Range("AW209").Value = 161.123079
$ws.Range("AW210").Value = 62.077384
$ws.Range("AK211").Value = 14.10985
$ws.Range("AK212").Value = 14.109838
$ws.Range("AQ213").Value = 19.061481
$ws.Range("AN214").Value = 6.107095
$ws.Range("AW215").Value = 69.33184
$ws.Range("AW216").Value = 121.097315
$ws.Range("AW217").Value = 138.338623
$ws.Range("AW218").Value = 78.10822899999999
$ws.Range("AW219").Value = 70.348426
$ws.Range("AW220").Value = 50.396319
$ws.Range("AW221").Value = 1.325185
$ws.Range("AW222").Value = 113.312975
$ws.Range("AW223").Value = 69.191435
$ws.Range("AK224").Value = 14.109988
$ws.Range("AK225").Value = 36.159502
$ws.Range("AQ226").Value = 26.064282
$ws.Range("AW227").Value = 103.33897
$ws.Range("AW228").Value = 149.01147
$ws.Range("AW229").Value = 114.005544
$ws.Range("AW230").Value = 78.039618
$ws.Range("AW231").Value = 21.294201
$ws.Range("AW232").Value = 50.099086
$ws.Range("AW233").Value = 44.378981
$ws.Range("AK234").Value = 14.109792
$ws.Range("AK235").Value = 36.959757
$ws.Range("AW236").Value = 114.005949
$ws.Range("AW237").Value = 85.266481
$ws.Range("AK238").Value = 14.109873
$ws.Range("AW239").Value = 134.159144
$ws.Range("AK240").Value = 14.10287
$ws.Range("AK241").Value = 14.109745
$ws.Range("AW242").Value = 113.273553
$ws.Range("AW243").Value = 44.388727
$ws.Range("AW244").Value = 44.386574
$ws.Range("AW245").Value = 18.971123
$ws.Range("AK246").Value = 14.109826
$ws.Range("AW247").Value = 111.217813
$ws.Range("AW248").Value = 35.377292
$ws.Range("AW249").Value = 21.149514
$ws.Range("AW250").Value = 142.13625
$ws.Range("AW251").Value = 152.084329
$ws.Range("AW252").Value = 131.000093
$ws.Range("AW253").Value = 21.150324
$ws.Range("AN254").Value = 7.275359
$ws.Range("AW255").Value = 7.158947
$ws.Range("AW256").Value = 131.000637
$ws.Range("AW257").Value = 50.390718
$ws.Range("AW258").Value = 2.137558
$ws.Range("AW259").Value = 141.029444
$ws.Range("AW260").Value = 56.244965
$ws.Range("AW261").Value = 8.331574
$ws.Range("AW262").Value = 127.221435
$ws.Range("AQ263").Value = 19.061389
$ws.Range("AW264").Value = 7.26463
$ws.Range("AQ265").Value = 1.324329
$ws.Range("AW266").Value = 36.380718
$ws.Range("AQ267").Value = 1.324213
$ws.Range("AW268").Value = 128.387708
$ws.Range("AW269").Value = 140.144375
$ws.Range("AW270").Value = 103.340035
$ws.Range("AW271").Value = 83.068102
$ws.Range("AW272").Value = 24.395521
$ws.Range("AW273").Value = 127.350104
$ws.Range("AK274").Value = 36.15941
$ws.Range("AW275").Value = 2.131215
$ws.Range("AW276").Value = 141.349306
$ws.Range("AW277").Value = 85.957199
$ws.Range("AW278").Value = 146.940208
$ws.Range("AW279").Value = 69.333056
$ws.Range("AW280").Value = 114.007882
$ws.Range("AK281").Value = 65.17569399999999
$ws.Range("AW282").Value = 138.339086
$ws.Range("AW283").Value = 69.325926
$ws.Range("AW284").Value = 44.137396
$ws.Range("AK285").Value = 14.11
$ws.Range("AW286").Value = 79.098866
$ws.Range("AW287").Value = 153.378808
$ws.Range("AW288").Value = 43.051019
$ws.Range("AK289").Value = 14.110069
$ws.Range("AK290").Value = 14.109757
$ws.Range("AQ291").Value = 0.061412
$ws.Range("AW292").Value = 43.049583
$ws.Range("AW293").Value = 43.051354
$ws.Range("AK294").Value = 36.959757
$ws.Range("AW295").Value = 146.94059
$ws.Range("AW296").Value = 100.00228
$ws.Range("AW297").Value = 19.08265
$ws.Range("AW298").Value = 100.005069
$ws.Range("AW299").Value = 113.310313
$ws.Range("AK300").Value = 9.13316
$ws.Range("AW301").Value = 156.273657
$ws.Range("AW302").Value = 50.389861
$ws.Range("AW303").Value = 14.099942
$ws.Range("AW304").Value = 2.135197
$ws.Range("AW305").Value = 135.047708
$ws.Range("AW306").Value = 90.106725
$ws.Range("AK307").Value = 14.109977
$ws.Range("AK308").Value = 14.109931
$ws.Range("AW309").Value = 114.007211
$ws.Range("AW310").Value = 83.01409700000001
$ws.Range("AW311").Value = 29.024815
$ws.Range("AW312").Value = 8.327788999999999
$ws.Range("AQ313").Value = 1.32853
$ws.Range("AW314").Value = 93.03169
$ws.Range("AK315").Value = 14.11
$ws.Range("AW316").Value = 85.263125
$ws.Range("AW317").Value = 19.082164
$ws.Range("AW318").Value = 161.132037
$ws.Range("AW319").Value = 71.331007
$ws.Range("AW320").Value = 50.397222
$ws.Range("AW321").Value = 0.045243
$ws.Range("AW322").Value = 118.370382
$ws.Range("AK323").Value = 14.109931
$ws.Range("AK324").Value = 14.109977
$ws.Range("AQ325").Value = 1.324468
$ws.Range("AW326").Value = 75.09076399999999
$ws.Range("AW327").Value = 54.029282
$ws.Range("AW328").Value = 19.083275
$ws.Range("AW329").Value = 0.04647
$ws.Range("AW330").Value = 127.350324
$ws.Range("AW331").Value = 56.244514
$ws.Range("AW332").Value = 44.378218
$ws.Range("AW333").Value = 50.450729
$ws.Range("AK334").Value = 14.110046
$ws.Range("AQ335").Value = 2.287882
$ws.Range("AW336").Value = 152.034942
$ws.Range("AW337").Value = 153.377118
$ws.Range("AK338").Value = 14.109942
$ws.Range("AK339").Value = 22.061377
$ws.Range("AW340").Value = 160.025787
$ws.Range("AW341").Value = 50.398438
$ws.Range("AW342").Value = 131.001157
$ws.Range("AW343").Value = 21.295706
$ws.Range("AW344").Value = 19.330741
$ws.Range("AW345").Value = 146.939942
$ws.Range("AW346").Value = 118.367697
$ws.Range("AW347").Value = 36.967894
$ws.Range("AW348").Value = 121.160255
$ws.Range("AW349").Value = 69.32728
$ws.Range("AW350").Value = 156.165694
$ws.Range("AW351").Value = 152.071551
$ws.Range("AW352").Value = 0.051863
$ws.Range("AW353").Value = 153.377917
$ws.Range("AW354").Value = 131.001678
$ws.Range("AW355").Value = 161.123438
$ws.Range("AW356").Value = 9.059815
$ws.Range("AW357").Value = 146.938472
$ws.Range("AW358").Value = 114.004352
$ws.Range("AW359").Value = 105.136979
$ws.Range("AW360").Value = 21.295023
$ws.Range("AK361").Value = 14.109965
$ws.Range("AK362").Value = 36.159491
$ws.Range("AW363").Value = 68.030509
$ws.Range("AW364").Value = 68.030046
$ws.Range("AK365").Value = 14.110081
$ws.Range("AQ366").Value = 19.061632
$ws.Range("AW367").Value = 0.045451
